# Trade #35 closed at 2026-02-17 15:23:11 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.73
$summary.Range("B4").Value = -0.27
$summary.Range("B5").Value = -0.15
$summary.Range("B6").Value = 35
$summary.Range("B7").Value = 10
$summary.Range("B9").Value = 28.57

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.73
$status.Range("D4").Value = 35
$status.Range("E4").Value = -0.27
$status.Range("F4").Value = -0.27
$status.Range("G4").Value = 28.57

# --- Append new trade row (#35) to both "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 36

    $ws.Cells.Item($row, 1).Value = 35

    # Keep date/time as plain text (matches existing inline-string formatting),
    # avoid Excel auto-converting "2026-02-17" into a date serial number.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"

    $ws.Cells.Item($row, 3).Value = "15:23:05"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.6899999999999999
    $ws.Cells.Item($row, 7).Value = 0.73
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 5.7971
    $ws.Cells.Item($row, 10).Value = 0.04
    $ws.Cells.Item($row, 11).Value = 99.73
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.15
}
